# Applies the NGE inference sheet edit:
#  1. Rename the worksheet tab.
#  2. Reorder header labels in columns G:K (left-rotate).
#  3. Fix a typo in the M1 header.
#  4. Left-rotate the numeric values in columns G:K for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "Inferencia por IA"

# 2. Reorder header row (G1:K1) - left rotate
$ws.Range("G1").Value = "qtd_carc"
$ws.Range("H1").Value = "qtd_tabela"
$ws.Range("I1").Value = "qtd_image"
$ws.Range("J1").Value = "qtd_estilos"
$ws.Range("K1").Value = "qtd_pag_word"

# 3. Fix typo in M1
$ws.Range("M1").Value = "QTD_PAG_INFERIDO"

# 4. Left-rotate values in G:K for each data row (rows 2-18)
for ($row = 2; $row -le 18; $row++) {
    $g = $ws.Cells.Item($row, 7).Value2
    $h = $ws.Cells.Item($row, 8).Value2
    $i = $ws.Cells.Item($row, 9).Value2
    $j = $ws.Cells.Item($row, 10).Value2
    $k = $ws.Cells.Item($row, 11).Value2

    $ws.Cells.Item($row, 7).Value = $h
    $ws.Cells.Item($row, 8).Value = $i
    $ws.Cells.Item($row, 9).Value = $j
    $ws.Cells.Item($row, 10).Value = $k
    $ws.Cells.Item($row, 11).Value = $g
}
